$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$ws.Range("B2:H2").Value = 0.051
$ws.Range("B3:H3").Value = 0.096
$ws.Range("B4:H4").Value = 0.045
$ws.Range("B6:H6").Value = 0.3
